$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 18946
$ws.Range("E3").Value = 5304
$ws.Range("E4").Value = 18016
$ws.Range("E5").Value = 1666
$ws.Range("E6").Value = 16533
$ws.Range("E7").Value = 3102
$ws.Range("E8").Value = 6841
$ws.Range("E9").Value = 1146
$ws.Range("E10").Value = 16797
$ws.Range("E11").Value = 15214
$ws.Range("E12").Value = 14293
$ws.Range("E13").Value = 12708
